# Update the NATMI Wnt4-Fzd6 LR-pair sheet following Dr Hou's advice:
# recompute rows 2-9 (new ligand/receptor-expressing-cell counts & derived
# stats) and append rows 10-13 for the previously-missing sender/target
# cluster combinations, growing the sheet from A1:T9 to A1:T13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate, G=Ligand average expression value,
# H=Ligand total expression value, I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate, M=Receptor average expression value,
# N=Receptor total expression value, O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity
$rows = @(
  ,@("ECs","Wnt4","Fzd6","ECs",3,1,1.437774666666667,4.313324,0.2965885866872326,0.2965885866872326,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,17.851346999248,160.662122993232,0.2550081218001355,0.2550081218001356)
  ,@("ECs","Wnt4","Fzd6","FAPs",3,1,1.437774666666667,4.313324,0.2965885866872326,0.2965885866872326,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,2.510181555781778,22.591634002036,0.03585817271627842,0.03585817271627843)
  ,@("ECs","Wnt4","Fzd6","sCs",3,1,1.437774666666667,4.313324,0.2965885866872326,0.2965885866872326,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.4005779206217777,3.605201285595999,0.005722292170818622,0.005722292170818623)
  ,@("FAPs","Wnt4","Fzd6","ECs",3,1,1.589574,4.768721999999999,0.3279022207198701,0.3279022207198702,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,19.736080842744,177.624727584696,0.2819317168399559,0.2819317168399559)
  ,@("FAPs","Wnt4","Fzd6","FAPs",3,1,1.589574,4.768721999999999,0.3279022207198701,0.3279022207198702,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,2.775204925261999,24.97684432735799,0.03964405574724195,0.03964405574724197)
  ,@("FAPs","Wnt4","Fzd6","sCs",3,1,1.589574,4.768721999999999,0.3279022207198701,0.3279022207198702,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.4428706822819999,3.985836140537999,0.006326448132672277,0.00632644813267228)
  ,@("M2","Wnt4","Fzd6","ECs",3,1,1.104440333333333,3.313321,0.2278273537140099,0.22782735371401,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,13.712682583292,123.414143249628,0.1958869227377649,0.1958869227377649)
  ,@("M2","Wnt4","Fzd6","FAPs",3,1,1.104440333333333,3.313321,0.2278273537140099,0.22782735371401,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,1.928219921013222,17.353979289119,0.02754479762764688,0.02754479762764689)
  ,@("M2","Wnt4","Fzd6","sCs",3,1,1.104440333333333,3.313321,0.2278273537140099,0.22782735371401,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.3077077531232222,2.769369778109,0.004395633348598187,0.004395633348598189)
  ,@("sCs","Wnt4","Fzd6","ECs",3,1,0.7159183333333333,2.147755,0.1476818388788872,0.1476818388788872,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,8.88881052626,79.99929473633999,0.1269774699597921,0.1269774699597921)
  ,@("sCs","Wnt4","Fzd6","FAPs",3,1,0.7159183333333333,2.147755,0.1476818388788872,0.1476818388788872,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,1.249907261160556,11.249165350445,0.0178550393483658,0.0178550393483658)
  ,@("sCs","Wnt4","Fzd6","sCs",3,1,0.7159183333333333,2.147755,0.1476818388788872,0.1476818388788872,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.1994617682105555,1.795155913895,0.002849329570729337,0.002849329570729339)
)

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r++
}
